$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 7: difference between the two angle readings for Exp 1 ---
$ws.Range("B7").Formula = "=B6-B5"

# --- New row 11: difference for Exp 2 + averaged k / dk stats (Exp1 & Exp2) ---
$ws.Range("B11").Formula = "=B10-B9"
$ws.Range("J11").Formula = "=AVERAGE(J10,J5)"
$ws.Range("K11").Formula = "=T.INV.2T(0.05,1)*STDEV.S(J5,J10)/SQRT(2)"
$ws.Range("M11").Formula = "=AVERAGE(M10,M5)"
$ws.Range("N11").Formula = "=T.INV.2T(0.05,1)*STDEV.S(M5,M10)/SQRT(2)"

# --- Brewster angle correction: use 57 degrees instead of 56.7, and fix the
#     uncertainty propagation formula for E13 (derivative of TAN wrt angle) ---
$ws.Range("B13").Formula = "=57*PI()/180"
$ws.Range("E13").Formula = "=C13/COS(B13)^2"

# --- Column K needs the same width/bestfit formatting as column H ---
$ws.Columns("K").ColumnWidth = $ws.Columns("H").ColumnWidth

# --- Update the view: drop the frozen/scrolled top-left cell and move the
#     active selection to G13 (Brewster angle area) instead of J35 ---
$ws.Range("G13").Select()
